$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.946.10'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '1.635.57'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = "'214.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D8').Value = "'28.83"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').Value = "'0.0904"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '1.867.52'
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('D13').Value = '1.630.61'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').Value = "'0.564"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = "'9.31"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +11.34%  '
$ws.Range('D16').Value = '29.962.23'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').Value = "'241.82"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = "'1.00"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').Value = "'4.15"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.19%  '
$ws.Range('D23').Value = "'9.84"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.32%  '
$ws.Range('E24').Value = '  +2.94%  '
$ws.Range('D25').Value = "'157.92"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').Value = "'15.49"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = "'6.58"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('D30').Value = "'0.0491"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.07%  '
$ws.Range('E31').Value = '  +4.03%  '
$ws.Range('E32').Value = '  +3.10%  '
$ws.Range('D34').Value = '1.430.14'
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('E35').Value = '  +4.72%  '
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('D37').Value = "'2.78"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.06%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'0.0171"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D39').Value = "'2.29"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').Value = "'75.65"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.02%  '
$ws.Range('D41').Value = "'0.553"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('D42').Value = "'1.99"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range('D43').Value = "'0.830"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').Value = "'0.0499"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('D47').Value = "'51.35"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.65%  '
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('D49').Value = '1.774.64'
$ws.Range('E49').Value = '  +1.65%  '
$ws.Range('E50').Value = '  +10.94%  '
$ws.Range('D51').Value = "'90.62"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.99%  '
